$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I0 in I1 and IF in J1, matching the style ---
# --- used by the existing header cells (e.g. H1 "IP").                 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data values for columns I and J, rows 2-23 ---
$iValues = @(3, 7, 4, 3, 3, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 5, 2, 8, 7, 9)
$jValues = @(5, 7, 8, 7, 9, 4, 7, 5, 4, 6, 4, 6, 7, 7, 6, 6, 6, 7, 5, 8, 7, 9)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
